$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "last row" border formatting (currently on row 21) onto row 18,
# since after the row deletion below, row 18 becomes the new last row.
$ws.Range("B21:J21").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)

# Update the "Valor Mora" total
$ws.Range("E11").Value = 170820

# Update "Cant. Trabajadores" and "Cant. Periodos"
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 1

# Row 16: first worker (replace Ana Francisca with Miledis)
$ws.Range("C16").Value = "45593147"
$ws.Range("D16").Value = "MILEDIS DOMINGUEZ BARRIOS"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 56940
$ws.Range("G16").Value = 1423500

# Row 17: second worker (new person)
$ws.Range("C17").Value = "1007257107"
$ws.Range("D17").Value = "DANIELA ISABEL CABALLERO ALVAREZ"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18: third worker (new person), already has row-21's formatting from above
$ws.Range("C18").Value = "1148440062"
$ws.Range("D18").Value = "YESICA BAENA TORRES"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 56940
$ws.Range("G18").Value = 1423500

# Remove the now-obsolete extra data rows (old rows 19,20,21)
$ws.Rows("19:21").Delete()
